$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HPPECbP")

# New row 7: "electrolysis with guaranteed clean electricity"
$ws.Range("A7").Value = "electrolysis with guaranteed clean electricity"
$ws.Range("A7").Font.Bold = $true

# New row 8: "natural gas reforming with CCS"
$ws.Range("A8").Value = "natural gas reforming with CCS"
$ws.Range("A8").Font.Bold = $true

# Every data cell in the two new rows repeats the same "excess capacity"
# assumption used by the rest of the table, so fill B7:AI8 with the same
# formula pattern as rows 2-6 (=$B$2).
$ws.Range("B7").Formula = "=`$B`$2"
$ws.Range("B8").Formula = "=`$B`$2"
$ws.Range("C7:AI8").Formula = "=`$B`$2"

# Restore the workbook's on-screen selection/scroll state for the HPPECbP
# sheet to cover the newly added rows.
[void]$ws.Range("B6:AI8").Select()

# Keep "About" as the active/selected tab, matching the saved workbook state.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
